$wb = $excel.ActiveWorkbook

# Hunk 0: ALC row 32
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2381.111
$ws.Range("I32").Value = 2075.4443
$ws.Range("J32").Value = 2686.7778
$ws.Range("K32").Value = 2075.4443
$ws.Range("L32").Value = 2686.7778
$ws.Range("M32").Value = -1749.4443
$ws.Range("N32").Value = -3338.7778

# Hunk 1: ALC row 52
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H52").Value = 7980.625
$ws.Range("I52").Value = 7977.857
$ws.Range("J52").Value = 8000
$ws.Range("K52").Value = 23933.571
$ws.Range("L52").Value = 24000
$ws.Range("M52").Value = -23773.571
$ws.Range("N52").Value = -24320

# Hunk 2: ALC row 61
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H61").Value = 1869
$ws.Range("I61").Value = 1869
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 5607
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -5435

# Hunk 3: ALC row 86
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 6866.222
$ws.Range("I86").Value = 6874.6665
$ws.Range("J86").Value = 6862
$ws.Range("K86").Value = 6874.6665
$ws.Range("L86").Value = 6862
$ws.Range("M86").Value = -5751.6665
$ws.Range("N86").Value = -9108

# Hunk 4: ALC row 89
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 6866.222
$ws.Range("I89").Value = 6874.6665
$ws.Range("J89").Value = 6862
$ws.Range("K89").Value = 34373.3325
$ws.Range("L89").Value = 34310
$ws.Range("M89").Value = -28757.3325
$ws.Range("N89").Value = -45542

# Hunk 5: ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 15156280
$ws.Range("I137").Value = 19236170
$ws.Range("J137").Value = 2408
$ws.Range("K137").Value = 57708510
$ws.Range("L137").Value = 7224
$ws.Range("M137").Value = -57705960

# Hunk 6: ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 5098.6865
$ws.Range("I138").Value = 2405.0715
$ws.Range("J138").Value = 7032.564
$ws.Range("K138").Value = 7215.2145
$ws.Range("L138").Value = 21097.692
$ws.Range("M138").Value = -2075.2145
$ws.Range("N138").Value = -31377.692

# Hunk 7: ARM row 23
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 1100001
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 1100001
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 1100001
$ws.Range("M23").ClearContents()
$ws.Range("N23").Value = -1100519

# Hunk 8: ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6585.829
$ws.Range("I32").Value = 5398.4863
$ws.Range("J32").Value = 17568.75
$ws.Range("K32").Value = 5398.4863
$ws.Range("L32").Value = 17568.75
$ws.Range("M32").Value = -5111.4863
$ws.Range("N32").Value = -18142.75

# Hunk 9: ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2731.476
$ws.Range("I61").Value = 2288.4375
$ws.Range("J61").Value = 4149.2
$ws.Range("K61").Value = 2288.4375
$ws.Range("L61").Value = 4149.2
$ws.Range("M61").Value = -2076.4375

# Hunk 10: ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 9860.833000000001
$ws.Range("I74").Value = 6268.3
$ws.Range("J74").Value = 27823.5
$ws.Range("K74").Value = 6268.3
$ws.Range("L74").Value = 27823.5
$ws.Range("M74").Value = -5394.3
$ws.Range("N74").Value = -29571.5

# Hunk 11: ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 9860.833000000001
$ws.Range("I77").Value = 6268.3
$ws.Range("J77").Value = 27823.5
$ws.Range("K77").Value = 31341.5
$ws.Range("L77").Value = 139117.5
$ws.Range("M77").Value = -26973.5
$ws.Range("N77").Value = -147853.5

# Hunk 12: ARM row 88
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 1442.7142
$ws.Range("I88").Value = 1249.5
$ws.Range("J88").Value = 1520
$ws.Range("K88").Value = 1249.5
$ws.Range("L88").Value = 1520
$ws.Range("M88").Value = -843.5
$ws.Range("N88").Value = -2332

# Hunk 13: ARM row 91
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 1442.7142
$ws.Range("I91").Value = 1249.5
$ws.Range("J91").Value = 1520
$ws.Range("K91").Value = 1249.5
$ws.Range("L91").Value = 1520
$ws.Range("M91").Value = 154.5
$ws.Range("N91").Value = -4328

# Hunk 14: ARM row 104
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H104").Value = 30000
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 30000
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 30000
$ws.Range("N104").Value = -36988

# Hunk 15: ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 14301.091
$ws.Range("I132").Value = 8078.346
$ws.Range("J132").Value = 37414.145
$ws.Range("K132").Value = 24235.038
$ws.Range("L132").Value = 112242.435
$ws.Range("M132").Value = -21705.038

# Hunk 16: ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2731.476
$ws.Range("I136").Value = 2288.4375
$ws.Range("J136").Value = 4149.2
$ws.Range("K136").Value = 6865.3125
$ws.Range("L136").Value = 12447.6
$ws.Range("M136").Value = -4315.3125

# Hunk 17: BSM row 26
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 3899.6667
$ws.Range("I26").Value = 3899.6667
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 3899.6667
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -3607.6667

# Hunk 18: BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4429.0303
$ws.Range("I86").Value = 4113.5415
$ws.Range("J86").Value = 5270.3335
$ws.Range("K86").Value = 4113.5415
$ws.Range("L86").Value = 5270.3335
$ws.Range("M86").Value = -2990.5415
$ws.Range("N86").Value = -7516.3335

# Hunk 19: BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 4429.0303
$ws.Range("I89").Value = 4113.5415
$ws.Range("J89").Value = 5270.3335
$ws.Range("K89").Value = 20567.7075
$ws.Range("L89").Value = 26351.6675
$ws.Range("M89").Value = -14951.7075
$ws.Range("N89").Value = -37583.6675

# Hunk 20: BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 55557996
$ws.Range("I134").Value = 71430570
$ws.Range("J134").Value = 3982
$ws.Range("K134").Value = 214291710
$ws.Range("L134").Value = 11946
$ws.Range("M134").Value = -214289175

# Hunk 21: CRP row 16
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2623.1765
$ws.Range("I16").Value = 1859.6
$ws.Range("J16").Value = 2941.3333
$ws.Range("K16").Value = 1859.6
$ws.Range("L16").Value = 2941.3333
$ws.Range("M16").Value = -1572.6

# Hunk 22: CRP row 22
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1451
$ws.Range("I22").Value = 176.83333
$ws.Range("J22").Value = 3999.3333
$ws.Range("K22").Value = 176.83333
$ws.Range("L22").Value = 3999.3333
$ws.Range("M22").Value = 173.16667
$ws.Range("N22").Value = -4699.3333

# Hunk 23: CRP row 86
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 76927576
$ws.Range("I86").Value = 142861040
$ws.Range("J86").Value = 5202.5
$ws.Range("K86").Value = 142861040
$ws.Range("L86").Value = 5202.5
$ws.Range("M86").Value = -142859917
$ws.Range("N86").Value = -7448.5

# Hunk 24: CRP row 89
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 76927576
$ws.Range("I89").Value = 142861040
$ws.Range("J89").Value = 5202.5
$ws.Range("K89").Value = 714305200
$ws.Range("L89").Value = 26012.5
$ws.Range("M89").Value = -714299584
$ws.Range("N89").Value = -37244.5

# Hunk 25: CRP row 113
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 2623.1765
$ws.Range("I113").Value = 1859.6
$ws.Range("J113").Value = 2941.3333
$ws.Range("K113").Value = 1859.6
$ws.Range("L113").Value = 2941.3333
$ws.Range("M113").Value = 310.4000000000001

# Hunk 26: CUL row 36
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 1000
$ws.Range("I36").Value = 1000
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 3000
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -2831

# Hunk 27: CUL row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 2687.818
$ws.Range("I122").Value = 5647.75
$ws.Range("J122").Value = 996.4286
$ws.Range("K122").Value = 50829.75
$ws.Range("L122").Value = 8967.857399999999
$ws.Range("M122").Value = -48379.75
$ws.Range("N122").Value = -13867.8574

# Hunk 28: GSM row 54
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H54").Value = 11000
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 11000
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 11000
$ws.Range("N54").Value = -11780

# Hunk 29: GSM row 55
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H55").Value = 16500
$ws.Range("I55").Value = 3000
$ws.Range("J55").Value = 30000
$ws.Range("K55").Value = 3000
$ws.Range("L55").Value = 30000
$ws.Range("M55").Value = -2673
$ws.Range("N55").Value = -30654

# Hunk 30: GSM row 117
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H117").Value = 43159.43
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 43159.43
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 43159.43
$ws.Range("N117").Value = -50043.43

# Hunk 31: GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 5164.8184
$ws.Range("I126").Value = 5217.1577
$ws.Range("J126").Value = 4833.3335
$ws.Range("K126").Value = 15651.4731
$ws.Range("L126").Value = 14500.0005
$ws.Range("M126").Value = -13181.4731
$ws.Range("N126").Value = -19440.0005

# Hunk 32: GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 5702.4
$ws.Range("I132").Value = 2430.1428
$ws.Range("J132").Value = 13337.667
$ws.Range("K132").Value = 7290.428400000001
$ws.Range("L132").Value = 40013.001
$ws.Range("M132").Value = -4760.428400000001

# Hunk 33: GSM row 139
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H139").Value = 99999.89999999999
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 99999.89999999999
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 99999.89999999999
$ws.Range("N139").Value = -110279.9

# Hunk 34: WVR row 107
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 823.13635
$ws.Range("I107").Value = 1062.3334
$ws.Range("J107").Value = 733.4375
$ws.Range("K107").Value = 3187.0002
$ws.Range("L107").Value = 2200.3125
$ws.Range("M107").Value = -1267.0002
$ws.Range("N107").Value = -6040.3125
